# Update "Shelter Assigned" (column C) and "Level" (column D) values
# for the allocation results sheet, per the re-run allocation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# r, Shelter Assigned, Level
$data = @(
    @(2,  "NV9 Multi-Purpose",               1),
    @(3,  "F. Mendoza Memorial Elem Sch.",    2),
    @(4,  "San Marcos National H.S.",         1),
    @(5,  "BMLTC Multi-Purpose Bldg and EC",  1),
    @(6,  "BMLTC Multi-Purpose Bldg and EC",  1),
    @(7,  "San Marcos National H.S.",         1),
    @(8,  "San Marcos National H.S.",         1),
    @(9,  "San Marcos National H.S.",         1),
    @(10, "F. Mendoza Memorial Elem Sch.",    2),
    @(11, "San Marcos National H.S.",         1),
    @(12, "Palimbang Empty Lot",              1),
    @(13, "Palimbang Empty Lot",              1),
    @(14, "San Marcos National H.S.",         1),
    @(15, "BMLTC Multi-Purpose Bldg and EC",  1),
    @(16, "San Marcos National H.S.",         1),
    @(17, "Mun. Covered Court",               1),
    @(18, "San Marcos Elem. Sch.",            1),
    @(19, "Mun. Covered Court",               1),
    @(20, "San Marcos National H.S.",         1),
    @(21, "BMLTC Multi-Purpose Bldg and EC",  1),
    @(22, "Palimbang Empty Lot",              1),
    @(23, "NV9 Multi-Purpose",                1),
    @(24, "San Marcos National H.S.",         1),
    @(25, "NV9 Multi-Purpose",                1),
    @(26, "F. Mendoza Memorial Elem Sch.",    2),
    @(27, "F. Mendoza Memorial Elem Sch.",    2),
    @(28, "F. Mendoza Memorial Elem Sch.",    2),
    @(29, "NV9 Multi-Purpose",                1),
    @(30, "F. Mendoza Memorial Elem Sch.",    2)
)

foreach ($row in $data) {
    $r = $row[0]
    $shelter = $row[1]
    $level = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $shelter
    $ws.Cells.Item($r, 4).Value2 = $level
}
